$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 107
$ws.Range("H107").Value = 935.8461
$ws.Range("I107").Value = 684.44446
$ws.Range("K107").Value = 684.44446
$ws.Range("M107").Value = 1235.55554
# row 132
$ws.Range("H132").Value = 1175.52
$ws.Range("I132").Value = 1175.52
$ws.Range("K132").Value = 3526.56
$ws.Range("M132").Value = -996.5599999999999
# row 136
$ws.Range("H136").Value = 65000
$ws.Range("J136").Value = 65000
$ws.Range("L136").Value = 65000
$ws.Range("N136").Value = -75200
# row 138
$ws.Range("H138").Value = 2005.75
$ws.Range("I138").Value = 1372.5238
$ws.Range("J138").Value = 3214.6365
$ws.Range("K138").Value = 4117.5714
$ws.Range("L138").Value = 9643.9095
$ws.Range("M138").Value = 1022.4286
$ws.Range("N138").Value = -19923.9095
# row 140
$ws.Range("H140").Value = 53871.855
$ws.Range("J140").Value = 53871.855
$ws.Range("L140").Value = 53871.855
$ws.Range("N140").Value = -64231.855
# row 141
$ws.Range("H141").Value = 3458.389
$ws.Range("I141").Value = 2422.818
$ws.Range("K141").Value = 7268.454000000001
$ws.Range("M141").Value = -2088.454000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 3883.625
$ws.Range("I32").Value = 2797.3057
$ws.Range("J32").Value = 7142.5835
$ws.Range("K32").Value = 2797.3057
$ws.Range("L32").Value = 7142.5835
$ws.Range("M32").Value = -2510.3057
$ws.Range("N32").Value = -7716.5835
# row 45
$ws.Range("H45").Value = 4501392.5
$ws.Range("I45").Value = 6429869.5
$ws.Range("J45").Value = 1613.8334
$ws.Range("K45").Value = 6429869.5
$ws.Range("L45").Value = 1613.8334
$ws.Range("M45").Value = -6429492.5
$ws.Range("N45").Value = -2367.8334
# row 88
$ws.Range("H88").Value = 2714.5264
$ws.Range("I88").Value = 2114.111
$ws.Range("J88").Value = 3254.9
$ws.Range("K88").Value = 2114.111
$ws.Range("L88").Value = 3254.9
$ws.Range("M88").Value = -1708.111
$ws.Range("N88").Value = -4066.9
# row 91
$ws.Range("H91").Value = 2714.5264
$ws.Range("I91").Value = 2114.111
$ws.Range("J91").Value = 3254.9
$ws.Range("K91").Value = 2114.111
$ws.Range("L91").Value = 3254.9
$ws.Range("M91").Value = -710.1109999999999
$ws.Range("N91").Value = -6062.9
# row 96
$ws.Range("H96").Value = 29967
$ws.Range("J96").Value = 29967
$ws.Range("L96").Value = 29967
$ws.Range("N96").Value = -35459
# row 110
$ws.Range("H110").Value = 2032.6666
$ws.Range("I110").Value = 650
$ws.Range("K110").Value = 650
$ws.Range("M110").Value = 1395
# row 119
$ws.Range("H119").Value = 44997.5
$ws.Range("J119").Value = 44997.5
$ws.Range("L119").Value = 44997.5
$ws.Range("N119").Value = -54673.5
# row 122
$ws.Range("H122").Value = 1850.8235
$ws.Range("I122").Value = 1891.5
$ws.Range("K122").Value = 5674.5
$ws.Range("M122").Value = -3224.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 184381.64
$ws.Range("I86").Value = 2525
$ws.Range("J86").Value = 669332.7
$ws.Range("K86").Value = 2525
$ws.Range("L86").Value = 669332.7
$ws.Range("M86").Value = -1402
$ws.Range("N86").Value = -671578.7
# row 89
$ws.Range("H89").Value = 184381.64
$ws.Range("I89").Value = 2525
$ws.Range("J89").Value = 669332.7
$ws.Range("K89").Value = 12625
$ws.Range("L89").Value = 3346663.5
$ws.Range("M89").Value = -7009
$ws.Range("N89").Value = -3357895.5
# row 107
$ws.Range("H107").Value = 1363.0526
$ws.Range("I107").Value = 1283.6875
$ws.Range("K107").Value = 1283.6875
$ws.Range("M107").Value = 636.3125
# row 134
$ws.Range("H134").Value = 9150
$ws.Range("I134").Value = 12081
$ws.Range("J134").Value = 3654.375
$ws.Range("K134").Value = 36243
$ws.Range("L134").Value = 10963.125
$ws.Range("M134").Value = -33708
$ws.Range("N134").Value = -16033.125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 5266.75
$ws.Range("I31").Value = 1906
$ws.Range("J31").Value = 5938.9
$ws.Range("K31").Value = 1906
$ws.Range("L31").Value = 5938.9
$ws.Range("M31").Value = -1611
$ws.Range("N31").Value = -6528.9
# row 34
$ws.Range("H34").Value = 5266.75
$ws.Range("I34").Value = 1906
$ws.Range("J34").Value = 5938.9
$ws.Range("K34").Value = 1906
$ws.Range("L34").Value = 5938.9
$ws.Range("M34").Value = -1704
$ws.Range("N34").Value = -6342.9
# row 93
$ws.Range("H93").Value = 6666.3335
$ws.Range("I93").Value = 4999.5
$ws.Range("K93").Value = 4999.5
$ws.Range("M93").Value = -3127.5
# row 99
$ws.Range("H99").Value = 1949.6666
$ws.Range("I99").Value = 1839.6
$ws.Range("K99").Value = 1839.6
$ws.Range("M99").Value = -341.5999999999999
# row 126
$ws.Range("H126").Value = 1949.6666
$ws.Range("I126").Value = 1839.6
$ws.Range("K126").Value = 5518.799999999999
$ws.Range("M126").Value = -3048.799999999999
# row 132
$ws.Range("H132").Value = 2680
$ws.Range("I132").Value = 1126.25
$ws.Range("K132").Value = 3378.75
$ws.Range("M132").Value = -848.75
# row 134
$ws.Range("H134").Value = 5319.6
$ws.Range("I134").Value = 4149.75
$ws.Range("K134").Value = 12449.25
$ws.Range("M134").Value = -9914.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 973.3333
$ws.Range("I5").Value = 585
$ws.Range("K5").Value = 1755
$ws.Range("M5").Value = -1643
# row 9
$ws.Range("H9").Value = 100000000
$ws.Range("J9").Value = 100000000
$ws.Range("L9").Value = 300000000
$ws.Range("N9").Value = -300000448
# row 55
$ws.Range("H55").Value = 35168
$ws.Range("I55").Value = 50252
$ws.Range("K55").Value = 150756
$ws.Range("M55").Value = -150579
# row 131
$ws.Range("H131").Value = 33380406
$ws.Range("J131").Value = 87862.25
$ws.Range("L131").Value = 263586.75
$ws.Range("N131").Value = -273666.75
# row 132
$ws.Range("H132").Value = 1911
$ws.Range("J132").Value = 2533
$ws.Range("L132").Value = 22797
$ws.Range("N132").Value = -27857
# row 135
$ws.Range("H135").Value = 973.3333
$ws.Range("I135").Value = 585
$ws.Range("K135").Value = 5265
$ws.Range("M135").Value = -2730
# row 137
$ws.Range("H137").Value = 3239.0356
$ws.Range("I137").Value = 1138.9166
$ws.Range("J137").Value = 4814.125
$ws.Range("K137").Value = 3416.7498
$ws.Range("L137").Value = 14442.375
$ws.Range("M137").Value = 1683.2502
$ws.Range("N137").Value = -24642.375

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 97
$ws.Range("H97").Value = 1194.4546
$ws.Range("I97").Value = 517.5
$ws.Range("J97").Value = 2999.6667
$ws.Range("K97").Value = 517.5
$ws.Range("L97").Value = 2999.6667
$ws.Range("M97").Value = -21.5
$ws.Range("N97").Value = -3991.6667
# row 110
$ws.Range("H110").Value = 99892
$ws.Range("J110").Value = 99892
$ws.Range("L110").Value = 99892
$ws.Range("N110").Value = -108072
# row 132
$ws.Range("H132").Value = 2914.04
$ws.Range("I132").Value = 2597.7693
$ws.Range("K132").Value = 7793.3079
$ws.Range("M132").Value = -5263.3079
# row 134
$ws.Range("H134").Value = 43123
$ws.Range("J134").Value = 43123
$ws.Range("L134").Value = 129369
$ws.Range("N134").Value = -134439
# row 141
$ws.Range("H141").Value = 22749.5
$ws.Range("J141").Value = 22749.5
$ws.Range("L141").Value = 22749.5
$ws.Range("N141").Value = -33109.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 3776.261
$ws.Range("I7").Value = 1920.4615
$ws.Range("J7").Value = 6188.8
$ws.Range("K7").Value = 1920.4615
$ws.Range("L7").Value = 6188.8
$ws.Range("M7").Value = -1808.4615
$ws.Range("N7").Value = -6412.8
# row 61
$ws.Range("H61").Value = 2343.8948
$ws.Range("I61").Value = 2234.3333
$ws.Range("J61").Value = 2531.7144
$ws.Range("K61").Value = 2234.3333
$ws.Range("L61").Value = 2531.7144
$ws.Range("M61").Value = -2032.3333
$ws.Range("N61").Value = -2935.7144
# row 113
$ws.Range("H113").Value = 2343.8948
$ws.Range("I113").Value = 2234.3333
$ws.Range("J113").Value = 2531.7144
$ws.Range("K113").Value = 2234.3333
$ws.Range("L113").Value = 2531.7144
$ws.Range("M113").Value = -64.33329999999978
$ws.Range("N113").Value = -6871.7144
# row 126
$ws.Range("H126").Value = 3776.261
$ws.Range("I126").Value = 1920.4615
$ws.Range("J126").Value = 6188.8
$ws.Range("K126").Value = 5761.3845
$ws.Range("L126").Value = 18566.4
$ws.Range("M126").Value = -3291.3845
$ws.Range("N126").Value = -23506.4
# row 132
$ws.Range("H132").Value = 1938.25
$ws.Range("I132").Value = 1874.25
$ws.Range("J132").Value = 1954.25
$ws.Range("K132").Value = 5622.75
$ws.Range("L132").Value = 5862.75
$ws.Range("M132").Value = -3092.75
$ws.Range("N132").Value = -10922.75
# row 136
$ws.Range("H136").Value = 5587.3
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 5587.3
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 16761.9
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -21861.9

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 19
$ws.Range("H19").Value = 16000
$ws.Range("J19").Value = 16000
$ws.Range("L19").Value = 16000
$ws.Range("N19").Value = -16348
# row 107
$ws.Range("H107").Value = 506.45456
$ws.Range("I107").Value = 423.21054
$ws.Range("K107").Value = 1269.63162
$ws.Range("M107").Value = 650.3683800000001
# row 126
$ws.Range("H126").Value = 5561.143
$ws.Range("I126").Value = 4799.125
$ws.Range("K126").Value = 14397.375
$ws.Range("M126").Value = -11927.375
# row 132
$ws.Range("H132").Value = 1587.1904
$ws.Range("I132").Value = 888.93335
$ws.Range("K132").Value = 2666.80005
$ws.Range("M132").Value = -136.8000499999998
# row 136
$ws.Range("H136").Value = 2401.6333
$ws.Range("I136").Value = 1676.8
$ws.Range("K136").Value = 5030.4
$ws.Range("M136").Value = -2480.4
# row 140
$ws.Range("H140").Value = 54999
$ws.Range("J140").Value = 54999
$ws.Range("L140").Value = 54999
$ws.Range("N140").Value = -65359
# row 141
$ws.Range("H141").Value = 70163.42999999999
$ws.Range("J141").Value = 70163.42999999999
$ws.Range("L141").Value = 70163.42999999999
$ws.Range("N141").Value = -80523.42999999999
